$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 (Objetivos:) - B/C text changes from the old (misplaced) docente name
# to the correct Portuguese objectives text. Row/label/style/height stay the same.
$ws.Range("B10:C10").Value = 'Propiciar ao aluno um panorama geral da área de Física do Estado Sólido, com ênfase nas idéias fundamentais e conceitos gerais, como gás de elétron, excitações elementares, estrutura de bandas, etc. O curso deve ser rico em resultados experimentais que ilustrem princípios e comportamentos gerais dos sólidos (por exemplo, comportamento das grandezas físicas com a temperatura).'

# Rows 13-26 are being substantially rebuilt/reshuffled, so start from a clean slate
# (this removes stray cells and old row heights for that block).
$ws.Range("A13:C26").Clear()

# --- Row 13-15: Docentes responsaveis (label already on row 12) ---
# Use row 2 (B/C-only data style) as the formatting template.
$ws.Range("B2:C2").Copy() | Out-Null
$ws.Range("B13:C15").PasteSpecial(-4122) | Out-Null
$ws.Range("B13:C13").Value = '5840730 - Antonio Jefferson da Silva Machado'
$ws.Range("B14:C14").Value = '5840726 - Cristina Bormio Nunes'
$ws.Range("B15:C15").Value = '1341653 - Maria José Ramos Sandim'
$ws.Rows(13).AutoFit()
$ws.Rows(14).AutoFit()
$ws.Rows(15).AutoFit()

# --- Row 16: Programa resumido (PT short syllabus) ---
# Use row 3 (A/B/C full data style) as the formatting template.
$ws.Range("A3:C3").Copy() | Out-Null
$ws.Range("A16:C16").PasteSpecial(-4122) | Out-Null
$ws.Range("A16").Value = "Programa resumido:"
$ws.Range("B16:C16").Value = 'Estrutura e ligações cristalinas. Vibrações da rede, fônons e propriedades térmicas. Gás de Fermi de elétrons livres. Bandas de energia. Semicondutores. Metais e superfícies de Fermi.'
$ws.Rows(16).RowHeight = 60

# --- Row 17: Short syllabus (EN) ---
$ws.Range("A3:C3").Copy() | Out-Null
$ws.Range("A17:C17").PasteSpecial(-4122) | Out-Null
$ws.Range("A17").Value = "Short syllabus:"
$ws.Range("B17:C17").Value = 'Crystal structure and bonds. Lattice vibrations, phonons and thermal properties. Free electron Fermi gas. Power bands. Semiconductors. Fermi metals and surfaces.'
$ws.Rows(17).RowHeight = 60

# --- Row 18: Programa (PT full syllabus) ---
$ws.Range("A3:C3").Copy() | Out-Null
$ws.Range("A18:C18").PasteSpecial(-4122) | Out-Null
$ws.Range("A18").Value = "Programa:"
$ws.Range("B18:C18").Value = '¨ Estrutura dos cristais.¨ Difração em cristais e a rede recíproca.¨ Ligações em cristais: cristais iônicos e cristais covalentes¨ Constantes elásticas e ondas elásticas.¨ Vibrações de cristais. Fônons¨ Gás de Fermi: modelo do elétron livre; movimento em campos magnéticos.¨ Bandas de energia. Funções de Bloch.¨ Cristais semicondutores.'
$ws.Rows(18).RowHeight = 120

# --- Row 19: Syllabus (EN full syllabus) ---
$ws.Range("A3:C3").Copy() | Out-Null
$ws.Range("A19:C19").PasteSpecial(-4122) | Out-Null
$ws.Range("A19").Value = "Syllabus:"
$ws.Range("B19:C19").Value = '¨ Structure of crystals.¨ Crystal diffraction and the reciprocal lattice.¨ Bonds in crystals: ionic crystals and covalent crystals¨ Elastic constants and elastic waves.¨ Crystal vibrations. phonons¨ Fermi gas: free electron model; movement in magnetic fields.¨ Energy bands. Bloch functions.¨ Semiconductor crystals.'
$ws.Rows(19).RowHeight = 120

# --- Row 20: Avaliacao (label only, A column) ---
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A20").PasteSpecial(-4122) | Out-Null
$ws.Range("A20").Value = "Avaliação:"
$ws.Rows(20).AutoFit()

# --- Row 21: Metodo ---
$ws.Range("A3:C3").Copy() | Out-Null
$ws.Range("A21:C21").PasteSpecial(-4122) | Out-Null
$ws.Range("A21").Value = "Método:"
$ws.Range("B21:C21").Value = 'Aulas expositivas, seminários e exercícios comentados.'
$ws.Rows(21).RowHeight = 60

# --- Row 22: Criterio ---
$ws.Range("A3:C3").Copy() | Out-Null
$ws.Range("A22:C22").PasteSpecial(-4122) | Out-Null
$ws.Range("A22").Value = "Critério:"
$ws.Range("B22:C22").Value = 'Média aritmética de duas provas com mesmo peso.'
$ws.Rows(22).RowHeight = 60

# --- Row 23: Norma de recuperacao ---
$ws.Range("A3:C3").Copy() | Out-Null
$ws.Range("A23:C23").PasteSpecial(-4122) | Out-Null
$ws.Range("A23").Value = "Norma de recuperação:"
$ws.Range("B23:C23").Value = 'Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação'
$ws.Rows(23).RowHeight = 60

# --- Row 24: Bibliografia ---
$ws.Range("A3:C3").Copy() | Out-Null
$ws.Range("A24:C24").PasteSpecial(-4122) | Out-Null
$ws.Range("A24").Value = "Bibliografia:"
$ws.Range("B24:C24").Value = 'ASHCROFT, N. W. Solid State Physics. Saunders College. KITTEL, C. Introduction to Solid State Physics. John Wiley & Sons. BLAKEMORE, J. S. Solid State Physics, Cambridge University Press. WERT,  C. A.; THOMSON, R. B. Physics of Solids. McGraw-Hill Book Co. Ltda. 1968. ZIMAN, J. M. Principles of the theory of solids, Cambridge, 2nd ed., 1972. SUTTON, A. P.  Electronic Structure of Materials, Oxford Science Publications.'
$ws.Rows(24).RowHeight = 120

# --- Row 25: Requisitos (label only) ---
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A25").PasteSpecial(-4122) | Out-Null
$ws.Range("A25").Value = "Requisitos:"
$ws.Rows(25).AutoFit()

# --- Row 26: Requisitos text (B/C only) ---
$ws.Range("B2:C2").Copy() | Out-Null
$ws.Range("B26:C26").PasteSpecial(-4122) | Out-Null
$ws.Range("B26:C26").Value = 'LOM3226 -  Mecânica Quântica  (Requisito fraco)
'
$ws.Rows(26).RowHeight = 30

$excel.CutCopyMode = 0

Write-Host "Edit complete"
